# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column H) and
# "Correspond Handback DateTime" (column K) values for the first
# localized file row (row 2) on both the "zh-cn" and "de-de" report
# sheets, reflecting a fresh handback report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 (08a4c95a-... file)
$wsZhCn.Range("H2").Value = "2016-08-23 12:48:35"
$wsZhCn.Range("K2").Value = "2016-08-23 12:48:52"

# de-de sheet: row 2 (08a4c95a-... file)
$wsDeDe.Range("H2").Value = "2016-08-23 12:48:40"
$wsDeDe.Range("K2").Value = "2016-08-23 12:49:00"
